$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# TTD sheet: separate the combined "Add/Edit" requirement notes in row 2
# into distinct "Add:" / "Edit:" lines (per commit message).
# ----------------------------------------------------------------------
$ttd = $wb.Worksheets.Item("TTD")

$nl = [char]10

$ttd.Range("A2").Value  = "Add: Required${nl}Edit: Required${nl}Edit Rates: Required${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional"
$ttd.Range("B2").Value  = "Add: Required${nl}Edit: Required${nl}Edit Rates: Optional${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional"
$ttd.Range("C2").Value  = "Add: Required${nl}Edit: Required${nl}Edit Rates: Optional${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional"
$ttd.Range("D2").Value  = "Add: Required${nl}Edit: Required${nl}Edit Rates: Optional${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional"
$ttd.Range("E2").Value  = "Add: Required${nl}Edit: Required${nl}Edit Rates: Optional${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional"
$ttd.Range("F2").Value  = "Add: Required${nl}Edit: Not Required${nl}Edit Rates: Required${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Required${nl}`"bombora`" or `"eyeota`" only"
$ttd.Range("G2").Value  = "Add: Required${nl}Edit: Not Required${nl}Edit Rates: Required${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Required"
$ttd.Range("H2").Value  = "Add: Required${nl}Edit: Not Required${nl}Edit Rates: Required${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional"
$ttd.Range("I2").Value  = "Add: Required${nl}Edit: Not Required${nl}Edit Rates: Required${nl}Retrieve Batch: Optional${nl}Retrieve Rates: Optional${nl}Values: CPM or PercentOfMediaCost"

# J2 gets a distinct format: no fill shading (previously shaded green), just the
# bordered / wrap-text / top-aligned look used elsewhere in the row.
$ttd.Range("K2").Copy() | Out-Null
$ttd.Range("J2").PasteSpecial(-4122) | Out-Null
$ttd.Range("J2").Value = "Add: Not Required${nl}Edit: Not Required${nl}Edit Rates: Optional${nl}Retrieve Batch: Required${nl}Retrieve Rates: Optional"
$ttd.Range("J2").WrapText = $true
$ttd.Range("J2").Font.Bold = $false
$ttd.Range("J2").Interior.Pattern = -4142

# Row 2 now wraps across more lines (long "Values: CPM or PercentOfMediaCost" line
# wraps in column I) -- bump the row height to match.
$ttd.Rows.Item(2).RowHeight = 119

# ----------------------------------------------------------------------
# Make TTD the active / selected sheet (it was "Adobe AAM" before).
# ----------------------------------------------------------------------
$ttd.Activate()
$ttd.Range("E2").Select() | Out-Null
